$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 125003460
$ws.Range("J64").Value = 4093.3333
$ws.Range("L64").Value = 4093.3333
$ws.Range("N64").Value = -4589.3333

# Row 67
$ws.Range("H67").Value = 125003460
$ws.Range("J67").Value = 4093.3333
$ws.Range("L67").Value = 4093.3333
$ws.Range("N67").Value = -5809.3333

# Row 74
$ws.Range("H74").Value = 3049.7646
$ws.Range("I74").Value = 2407.6428
$ws.Range("K74").Value = 2407.6428
$ws.Range("M74").Value = -1471.6428

# Row 77
$ws.Range("H77").Value = 3049.7646
$ws.Range("I77").Value = 2407.6428
$ws.Range("K77").Value = 12038.214
$ws.Range("M77").Value = -7358.214

# Row 138
$ws.Range("H138").Value = 1649.28
$ws.Range("I138").Value = 890.53625
$ws.Range("J138").Value = 3338.0967
$ws.Range("K138").Value = 2671.60875
$ws.Range("L138").Value = 10014.2901
$ws.Range("M138").Value = 2468.39125
$ws.Range("N138").Value = -20294.2901

# Row 141
$ws.Range("H141").Value = 5666.5
$ws.Range("I141").Value = 3131.0571
$ws.Range("J141").Value = 35246.668
$ws.Range("K141").Value = 9393.1713
$ws.Range("L141").Value = 105740.004
$ws.Range("M141").Value = -4213.1713
$ws.Range("N141").Value = -116100.004


# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 2109.8462
$ws.Range("I88").Value = 1454.2
$ws.Range("J88").Value = 2519.625
$ws.Range("K88").Value = 1454.2
$ws.Range("L88").Value = 2519.625
$ws.Range("M88").Value = -1048.2
$ws.Range("N88").Value = -3331.625

# Row 91
$ws.Range("H91").Value = 2109.8462
$ws.Range("I91").Value = 1454.2
$ws.Range("J91").Value = 2519.625
$ws.Range("K91").Value = 1454.2
$ws.Range("L91").Value = 2519.625
$ws.Range("M91").Value = -50.20000000000005
$ws.Range("N91").Value = -5327.625


# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 2718.6738
$ws.Range("I94").Value = 476.3243
$ws.Range("J94").Value = 11937.223
$ws.Range("K94").Value = 476.3243
$ws.Range("L94").Value = 11937.223
$ws.Range("M94").Value = -25.32429999999999
$ws.Range("N94").Value = -12839.223


# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 14304972
$ws.Range("I31").Value = 58824800
$ws.Range("K31").Value = 58824800
$ws.Range("M31").Value = -58824505

# Row 34
$ws.Range("H34").Value = 14304972
$ws.Range("I34").Value = 58824800
$ws.Range("K34").Value = 58824800
$ws.Range("M34").Value = -58824598

# Row 58
$ws.Range("H58").Value = 21277436
$ws.Range("I58").Value = 37037620
$ws.Range("J58").Value = 1184.2
$ws.Range("K58").Value = 37037620
$ws.Range("L58").Value = 1184.2
$ws.Range("M58").Value = -37037417
$ws.Range("N58").Value = -1590.2

# Row 136
$ws.Range("H136").Value = 21277436
$ws.Range("I136").Value = 37037620
$ws.Range("J136").Value = 1184.2
$ws.Range("K136").Value = 111112860
$ws.Range("L136").Value = 3552.6
$ws.Range("M136").Value = -111110310
$ws.Range("N136").Value = -8652.6


# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 693.6585
$ws.Range("I5").Value = 417.03226
$ws.Range("J5").Value = 1551.2
$ws.Range("K5").Value = 1251.09678
$ws.Range("L5").Value = 4653.6
$ws.Range("M5").Value = -1139.09678
$ws.Range("N5").Value = -4877.6

# Row 39
$ws.Range("H39").Value = 2030
$ws.Range("J39").Value = 2375
$ws.Range("L39").Value = 7125
$ws.Range("N39").Value = -7713

# Row 120
$ws.Range("H120").Value = 7161
$ws.Range("I120").Value = 5593.2
$ws.Range("J120").Value = 15000
$ws.Range("K120").Value = 16779.6
$ws.Range("L120").Value = 45000
$ws.Range("M120").Value = -11941.6
$ws.Range("N120").Value = -54676

# Row 135
$ws.Range("H135").Value = 693.6585
$ws.Range("I135").Value = 417.03226
$ws.Range("J135").Value = 1551.2
$ws.Range("K135").Value = 3753.29034
$ws.Range("L135").Value = 13960.8
$ws.Range("M135").Value = -1218.29034
$ws.Range("N135").Value = -19030.8


# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4234.2666
$ws.Range("I70").Value = 4064.25
$ws.Range("J70").Value = 4428.5713
$ws.Range("K70").Value = 4064.25
$ws.Range("L70").Value = 4428.5713
$ws.Range("M70").Value = -3794.25
$ws.Range("N70").Value = -4968.5713

# Row 73
$ws.Range("H73").Value = 4234.2666
$ws.Range("I73").Value = 4064.25
$ws.Range("J73").Value = 4428.5713
$ws.Range("K73").Value = 4064.25
$ws.Range("L73").Value = 4428.5713
$ws.Range("M73").Value = -3128.25
$ws.Range("N73").Value = -6300.5713

# Row 80
$ws.Range("H80").Value = 2571.3
$ws.Range("I80").Value = 2202.6
$ws.Range("J80").Value = 2940
$ws.Range("K80").Value = 2202.6
$ws.Range("L80").Value = 2940
$ws.Range("M80").Value = -1204.6
$ws.Range("N80").Value = -4936

# Row 83
$ws.Range("H83").Value = 2571.3
$ws.Range("I83").Value = 2202.6
$ws.Range("J83").Value = 2940
$ws.Range("K83").Value = 11013
$ws.Range("L83").Value = 14700
$ws.Range("M83").Value = -6021
$ws.Range("N83").Value = -24684


# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 305
$ws.Range("I22").Value = 210
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 210
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = 85
$ws.Range("N22").Value = -990

# Row 27
$ws.Range("H27").Value = 305
$ws.Range("I27").Value = 210
$ws.Range("J27").Value = 400
$ws.Range("K27").Value = 210
$ws.Range("L27").Value = 400
$ws.Range("M27").Value = -103
$ws.Range("N27").Value = -614

# Row 46
$ws.Range("H46").Value = 1381.0454
$ws.Range("I46").Value = 1155.7858
$ws.Range("J46").Value = 1775.25
$ws.Range("K46").Value = 1155.7858
$ws.Range("L46").Value = 1775.25
$ws.Range("M46").Value = -967.7858000000001
$ws.Range("N46").Value = -2151.25

# Row 81
$ws.Range("H81").Value = 41836.2
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 41836.2
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 41836.2
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -43832.2

# Row 82
$ws.Range("H82").Value = 2240.4
$ws.Range("I82").Value = 1451
$ws.Range("K82").Value = 1451
$ws.Range("M82").Value = -1090

# Row 84
$ws.Range("H84").Value = 41836.2
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 41836.2
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 125508.6
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -135492.6

# Row 85
$ws.Range("H85").Value = 2240.4
$ws.Range("I85").Value = 1451
$ws.Range("K85").Value = 1451
$ws.Range("M85").Value = -203


# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 47
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

# Row 132
$ws.Range("H132").Value = 3765.2983
$ws.Range("I132").Value = 4451.725
$ws.Range("J132").Value = 2150.1765
$ws.Range("K132").Value = 13355.175
$ws.Range("L132").Value = 6450.529500000001
$ws.Range("M132").Value = -10825.175
$ws.Range("N132").Value = -11510.5295

